$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MCH171-1"
$ws.Range("C2").Value = "PUBLICATIONS, NEWSLETTERS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 22C | GRAP COUNT NUMER: NONE"
